$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "In-house data 2023"
$ws.Range("A13").Value = "In-house data 2023"
$ws.Range("A4").Value = "Papadopoulos et al., 2012"
